# Auto-generated PowerShell COM-interop script to apply market-data cell updates
# across the Exodus_Profits workbook sheets (scheduled-runner style refresh).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 261.53845  # H28
$ws.Cells.Item(28, 9).Value = 199.5  # I28
$ws.Cells.Item(28, 10).Value = 1006  # J28
$ws.Cells.Item(28, 11).Value = 199.5  # K28
$ws.Cells.Item(28, 12).Value = 1006  # L28
$ws.Cells.Item(28, 13).Value = 285.5  # M28
$ws.Cells.Item(28, 14).Value = -1976  # N28
$ws.Cells.Item(76, 8).Value = 5300.387  # H76
$ws.Cells.Item(76, 9).Value = 5300.387  # I76
$ws.Cells.Item(76, 11).Value = 5300.387  # K76
$ws.Cells.Item(76, 13).Value = -4985.387  # M76
$ws.Cells.Item(79, 8).Value = 5300.387  # H79
$ws.Cells.Item(79, 9).Value = 5300.387  # I79
$ws.Cells.Item(79, 11).Value = 5300.387  # K79
$ws.Cells.Item(79, 13).Value = -4208.387  # M79
$ws.Cells.Item(107, 8).Value = 77729.46000000001  # H107
$ws.Cells.Item(107, 9).Value = 143541.28  # I107
$ws.Cells.Item(107, 10).Value = 949  # J107
$ws.Cells.Item(107, 11).Value = 143541.28  # K107
$ws.Cells.Item(107, 12).Value = 949  # L107
$ws.Cells.Item(107, 13).Value = -141621.28  # M107
$ws.Cells.Item(107, 14).Value = -4789  # N107
$ws.Cells.Item(132, 8).Value = 1220.6  # H132
$ws.Cells.Item(132, 9).Value = 968.7234  # I132
$ws.Cells.Item(132, 11).Value = 2906.1702  # K132
$ws.Cells.Item(132, 13).Value = -376.1702  # M132
$ws.Cells.Item(133, 8).Value = 92509.62  # H133
$ws.Cells.Item(133, 10).Value = 92509.62  # J133
$ws.Cells.Item(133, 12).Value = 92509.62  # L133
$ws.Cells.Item(133, 14).Value = -102629.62  # N133
$ws.Cells.Item(134, 8).Value = 99990  # H134
$ws.Cells.Item(134, 10).Value = 99990  # J134
$ws.Cells.Item(134, 12).Value = 99990  # L134
$ws.Cells.Item(134, 14).Value = -110130  # N134
$ws.Cells.Item(136, 8).Value = 67959.25  # H136
$ws.Cells.Item(136, 10).Value = 73382  # J136
$ws.Cells.Item(136, 12).Value = 73382  # L136
$ws.Cells.Item(136, 14).Value = -83582  # N136
$ws.Cells.Item(137, 8).Value = 691927.9  # H137
$ws.Cells.Item(137, 10).Value = 854432.3  # J137
$ws.Cells.Item(137, 12).Value = 2563296.9  # L137
$ws.Cells.Item(137, 14).Value = -2568396.9  # N137
$ws.Cells.Item(138, 8).Value = 10640059  # H138
$ws.Cells.Item(138, 9).Value = 1416.6428  # I138
$ws.Cells.Item(138, 10).Value = 12501821  # J138
$ws.Cells.Item(138, 11).Value = 4249.928400000001  # K138
$ws.Cells.Item(138, 12).Value = 37505463  # L138
$ws.Cells.Item(138, 13).Value = 890.0715999999993  # M138
$ws.Cells.Item(138, 14).Value = -37515743  # N138
$ws.Cells.Item(139, 8).Value = 53247  # H139
$ws.Cells.Item(139, 10).Value = 53247  # J139
$ws.Cells.Item(139, 12).Value = 53247  # L139
$ws.Cells.Item(139, 14).Value = -63527  # N139
$ws.Cells.Item(140, 8).Value = 91135.71000000001  # H140
$ws.Cells.Item(140, 10).Value = 91135.71000000001  # J140
$ws.Cells.Item(140, 12).Value = 91135.71000000001  # L140
$ws.Cells.Item(140, 14).Value = -101495.71  # N140
$ws.Cells.Item(141, 8).Value = 4867.5  # H141
$ws.Cells.Item(141, 9).Value = 4246.5713  # I141
$ws.Cells.Item(141, 10).Value = 6316.3335  # J141
$ws.Cells.Item(141, 11).Value = 12739.7139  # K141
$ws.Cells.Item(141, 12).Value = 18949.0005  # L141
$ws.Cells.Item(141, 13).Value = -7559.713899999999  # M141
$ws.Cells.Item(141, 14).Value = -29309.0005  # N141

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 13194.111  # H45
$ws.Cells.Item(45, 9).Value = 15535.286  # I45
$ws.Cells.Item(45, 11).Value = 15535.286  # K45
$ws.Cells.Item(45, 13).Value = -15158.286  # M45
$ws.Cells.Item(110, 8).Value = 1141.1428  # H110
$ws.Cells.Item(110, 9).Value = 1098.1538  # I110
$ws.Cells.Item(110, 11).Value = 1098.1538  # K110
$ws.Cells.Item(110, 13).Value = 946.8462  # M110
$ws.Cells.Item(132, 8).Value = 2437.5151  # H132
$ws.Cells.Item(132, 9).Value = 2073.375  # I132
$ws.Cells.Item(132, 10).Value = 3408.5557  # J132
$ws.Cells.Item(132, 11).Value = 6220.125  # K132
$ws.Cells.Item(132, 12).Value = 10225.6671  # L132
$ws.Cells.Item(132, 13).Value = -3690.125  # M132
$ws.Cells.Item(132, 14).Value = -15285.6671  # N132

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1117639.4  # H99
$ws.Cells.Item(99, 9).Value = 1047.3334  # I99
$ws.Cells.Item(99, 10).Value = 3127505  # J99
$ws.Cells.Item(99, 11).Value = 1047.3334  # K99
$ws.Cells.Item(99, 12).Value = 3127505  # L99
$ws.Cells.Item(99, 13).Value = 450.6666  # M99
$ws.Cells.Item(99, 14).Value = -3130501  # N99
$ws.Cells.Item(132, 8).Value = 58730  # H132
$ws.Cells.Item(132, 10).Value = 58730  # J132
$ws.Cells.Item(132, 12).Value = 58730  # L132
$ws.Cells.Item(132, 14).Value = -68850  # N132
$ws.Cells.Item(134, 8).Value = 2195.125  # H134
$ws.Cells.Item(134, 9).Value = 1754.7778  # I134
$ws.Cells.Item(134, 10).Value = 3516.1667  # J134
$ws.Cells.Item(134, 11).Value = 5264.3334  # K134
$ws.Cells.Item(134, 12).Value = 10548.5001  # L134
$ws.Cells.Item(134, 13).Value = -2729.3334  # M134
$ws.Cells.Item(134, 14).Value = -15618.5001  # N134
$ws.Cells.Item(135, 8).Value = 69054.5  # H135
$ws.Cells.Item(135, 10).Value = 69054.5  # J135
$ws.Cells.Item(135, 12).Value = 69054.5  # L135
$ws.Cells.Item(135, 14).Value = -79194.5  # N135
$ws.Cells.Item(138, 8).Value = 71825.5  # H138
$ws.Cells.Item(138, 10).Value = 71825.5  # J138
$ws.Cells.Item(138, 12).Value = 71825.5  # L138
$ws.Cells.Item(138, 14).Value = -82105.5  # N138
$ws.Cells.Item(140, 8).Value = 114202.664  # H140
$ws.Cells.Item(140, 10).Value = 65978  # J140
$ws.Cells.Item(140, 12).Value = 65978  # L140
$ws.Cells.Item(140, 14).Value = -76338  # N140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 49999.5  # H23
$ws.Cells.Item(23, 9).Value = 50000  # I23
$ws.Cells.Item(23, 10).Value = 49999  # J23
$ws.Cells.Item(23, 11).Value = 50000  # K23
$ws.Cells.Item(23, 12).Value = 49999  # L23
$ws.Cells.Item(23, 13).Value = -49760  # M23
$ws.Cells.Item(23, 14).Value = -50479  # N23
$ws.Cells.Item(27, 8).Value = 49999.5  # H27
$ws.Cells.Item(27, 9).Value = 50000  # I27
$ws.Cells.Item(27, 10).Value = 49999  # J27
$ws.Cells.Item(27, 11).Value = 50000  # K27
$ws.Cells.Item(27, 12).Value = 49999  # L27
$ws.Cells.Item(27, 13).Value = -49808  # M27
$ws.Cells.Item(27, 14).Value = -50383  # N27
$ws.Cells.Item(58, 8).Value = 2123.724  # H58
$ws.Cells.Item(58, 9).Value = 1881.5714  # I58
$ws.Cells.Item(58, 10).Value = 2349.7334  # J58
$ws.Cells.Item(58, 11).Value = 1881.5714  # K58
$ws.Cells.Item(58, 12).Value = 2349.7334  # L58
$ws.Cells.Item(58, 13).Value = -1678.5714  # M58
$ws.Cells.Item(58, 14).Value = -2755.7334  # N58
$ws.Cells.Item(132, 8).Value = 2363.0588  # H132
$ws.Cells.Item(132, 9).Value = 1166.9166  # I132
$ws.Cells.Item(132, 10).Value = 5233.8  # J132
$ws.Cells.Item(132, 11).Value = 3500.7498  # K132
$ws.Cells.Item(132, 12).Value = 15701.4  # L132
$ws.Cells.Item(132, 13).Value = -970.7498000000001  # M132
$ws.Cells.Item(132, 14).Value = -20761.4  # N132
$ws.Cells.Item(134, 8).Value = 1349.2881  # H134
$ws.Cells.Item(134, 9).Value = 815  # I134
$ws.Cells.Item(134, 11).Value = 2445  # K134
$ws.Cells.Item(134, 13).Value = 90  # M134
$ws.Cells.Item(136, 8).Value = 2123.724  # H136
$ws.Cells.Item(136, 9).Value = 1881.5714  # I136
$ws.Cells.Item(136, 10).Value = 2349.7334  # J136
$ws.Cells.Item(136, 11).Value = 5644.7142  # K136
$ws.Cells.Item(136, 12).Value = 7049.2002  # L136
$ws.Cells.Item(136, 13).Value = -3094.7142  # M136
$ws.Cells.Item(136, 14).Value = -12149.2002  # N136
$ws.Cells.Item(138, 8).Value = 58533.9  # H138
$ws.Cells.Item(138, 10).Value = 59704.332  # J138
$ws.Cells.Item(138, 12).Value = 59704.332  # L138
$ws.Cells.Item(138, 14).Value = -69984.33199999999  # N138

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 21047656  # H4
$ws.Cells.Item(4, 10).Value = 628.5714  # J4
$ws.Cells.Item(4, 12).Value = 1885.7142  # L4
$ws.Cells.Item(4, 14).Value = -2109.7142  # N4
$ws.Cells.Item(132, 8).Value = 2808.0232  # H132
$ws.Cells.Item(132, 9).Value = 2160.4443  # I132
$ws.Cells.Item(132, 10).Value = 2979.4412  # J132
$ws.Cells.Item(132, 11).Value = 19443.9987  # K132
$ws.Cells.Item(132, 12).Value = 26814.9708  # L132
$ws.Cells.Item(132, 13).Value = -16913.9987  # M132
$ws.Cells.Item(132, 14).Value = -31874.9708  # N132

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 1289151.9  # H11
$ws.Cells.Item(11, 9).Value = 380042.06  # I11
$ws.Cells.Item(11, 10).Value = 4167999.8  # J11
$ws.Cells.Item(11, 11).Value = 380042.06  # K11
$ws.Cells.Item(11, 12).Value = 4167999.8  # L11
$ws.Cells.Item(11, 13).Value = -379903.06  # M11
$ws.Cells.Item(11, 14).Value = -4168277.8  # N11
$ws.Cells.Item(12, 8).Value = 3507.1667  # H12
$ws.Cells.Item(12, 9).Value = 2766.3333  # I12
$ws.Cells.Item(12, 11).Value = 2766.3333  # K12
$ws.Cells.Item(12, 13).Value = -2626.3333  # M12
$ws.Cells.Item(14, 8).Value = 500000  # H14
$ws.Cells.Item(14, 9).Value = 500000  # I14
$ws.Cells.Item(14, 11).Value = 500000  # K14
$ws.Cells.Item(14, 13).Value = -499832  # M14
$ws.Cells.Item(29, 8).Value = 6000  # H29
$ws.Cells.Item(29, 9).Value = 6000  # I29
$ws.Cells.Item(29, 10).Value = 0  # J29
$ws.Cells.Item(29, 11).Value = 6000  # K29
$ws.Cells.Item(29, 12).Value = 0  # L29
$ws.Cells.Item(29, 13).Value = -5710  # M29
$ws.Cells.Item(29, 14).ClearContents()  # N29
$ws.Cells.Item(49, 8).Value = 55555  # H49
$ws.Cells.Item(49, 9).Value = 0  # I49
$ws.Cells.Item(49, 11).Value = 0  # K49
$ws.Cells.Item(49, 13).ClearContents()  # M49
$ws.Cells.Item(80, 8).Value = 3168.5  # H80
$ws.Cells.Item(80, 9).Value = 3183.8572  # I80
$ws.Cells.Item(80, 11).Value = 3183.8572  # K80
$ws.Cells.Item(80, 13).Value = -2185.8572  # M80
$ws.Cells.Item(83, 8).Value = 3168.5  # H83
$ws.Cells.Item(83, 9).Value = 3183.8572  # I83
$ws.Cells.Item(83, 11).Value = 15919.286  # K83
$ws.Cells.Item(83, 13).Value = -10927.286  # M83
$ws.Cells.Item(109, 8).Value = 57379.848  # H109
$ws.Cells.Item(109, 10).Value = 57379.848  # J109
$ws.Cells.Item(109, 12).Value = 57379.848  # L109
$ws.Cells.Item(109, 14).Value = -59459.848  # N109
$ws.Cells.Item(132, 8).Value = 1601.579  # H132
$ws.Cells.Item(132, 9).Value = 1215.25  # I132
$ws.Cells.Item(132, 11).Value = 3645.75  # K132
$ws.Cells.Item(132, 13).Value = -1115.75  # M132
$ws.Cells.Item(140, 8).Value = 65217.777  # H140
$ws.Cells.Item(140, 10).Value = 73326.664  # J140
$ws.Cells.Item(140, 12).Value = 73326.664  # L140
$ws.Cells.Item(140, 14).Value = -83686.664  # N140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1368.1904  # H46
$ws.Cells.Item(46, 9).Value = 1301.7142  # I46
$ws.Cells.Item(46, 10).Value = 1401.4286  # J46
$ws.Cells.Item(46, 11).Value = 1301.7142  # K46
$ws.Cells.Item(46, 12).Value = 1401.4286  # L46
$ws.Cells.Item(46, 13).Value = -1113.7142  # M46
$ws.Cells.Item(46, 14).Value = -1777.4286  # N46
$ws.Cells.Item(50, 8).Value = 9000  # H50
$ws.Cells.Item(50, 10).Value = 9000  # J50
$ws.Cells.Item(50, 12).Value = 9000  # L50
$ws.Cells.Item(50, 14).Value = -10274  # N50
$ws.Cells.Item(82, 8).Value = 2067.111  # H82
$ws.Cells.Item(82, 9).Value = 1568.875  # I82
$ws.Cells.Item(82, 10).Value = 2465.7  # J82
$ws.Cells.Item(82, 11).Value = 1568.875  # K82
$ws.Cells.Item(82, 12).Value = 2465.7  # L82
$ws.Cells.Item(82, 13).Value = -1207.875  # M82
$ws.Cells.Item(82, 14).Value = -3187.7  # N82
$ws.Cells.Item(85, 8).Value = 2067.111  # H85
$ws.Cells.Item(85, 9).Value = 1568.875  # I85
$ws.Cells.Item(85, 10).Value = 2465.7  # J85
$ws.Cells.Item(85, 11).Value = 1568.875  # K85
$ws.Cells.Item(85, 12).Value = 2465.7  # L85
$ws.Cells.Item(85, 13).Value = -320.875  # M85
$ws.Cells.Item(85, 14).Value = -4961.7  # N85
$ws.Cells.Item(122, 8).Value = 6076897.5  # H122
$ws.Cells.Item(122, 10).Value = 22226564  # J122
$ws.Cells.Item(122, 12).Value = 66679692  # L122
$ws.Cells.Item(122, 14).Value = -66684592  # N122
$ws.Cells.Item(132, 8).Value = 2433.9395  # H132
$ws.Cells.Item(132, 10).Value = 3344.9167  # J132
$ws.Cells.Item(132, 12).Value = 10034.7501  # L132
$ws.Cells.Item(132, 14).Value = -15094.7501  # N132

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 8751.5  # H9
$ws.Cells.Item(9, 9).Value = 1666.3334  # I9
$ws.Cells.Item(9, 10).Value = 30007  # J9
$ws.Cells.Item(9, 11).Value = 1666.3334  # K9
$ws.Cells.Item(9, 12).Value = 30007  # L9
$ws.Cells.Item(9, 13).Value = -1526.3334  # M9
$ws.Cells.Item(9, 14).Value = -30287  # N9
$ws.Cells.Item(44, 8).Value = 40499.5  # H44
$ws.Cells.Item(44, 10).Value = 40499.5  # J44
$ws.Cells.Item(44, 12).Value = 40499.5  # L44
$ws.Cells.Item(44, 14).Value = -41607.5  # N44
